# Apply crypto price/volume updates per commit:
# "Updated cryptos list on Wed Sep 18 08:51:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.147.33'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').Value = '2.319.22'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''543.69'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('D9').Value = '2.316.97'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').Value = '2.734.65'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').Value = '60.126.69'
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '2.309.13'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').Value = '''10.53'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').Value = '''313.45'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '''63.75'
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('E27').Value = '  -2.40%  '
$ws.Range('D28').Value = '''1.35'
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('D29').Value = '''1.20'
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').Value = '''172.04'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = '0.0₃0730'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('D33').Value = '''5.88'
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('E34').Value = '  +5.91%  '
$ws.Range('E35').Value = '  -1.99%  '
$ws.Range('D37').Value = '''17.72'
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('D40').Value = '''317.44'
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''1.52'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '''37.84'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('D43').Value = '''137.21'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '''0.0940'
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('D46').Value = '''18.92'
$ws.Range('E46').Value = '  +2.76%  '
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '''0.0493'
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0231'
$ws.Range('E49').Value = '  +23.12%  '
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('E51').Value = '  +0.32%  '
